# Auto update Excel log
# Appends new sensor log rows to the mmWave(InBed), mmWave(BR) and mmWave(HR)
# worksheets, mirroring the format of the existing rows (text for the
# Date/Timestamp/Hour/Location/Status columns, and for mmWave(InBed) the
# Value column too; numeric for the Value column on mmWave(BR)/mmWave(HR)).

$wb = $excel.ActiveWorkbook

function Add-LogRow($ws, $row, $date, $timestamp, $hour, $location, $value, $valueIsText, $status) {
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $date

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $timestamp

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $hour

    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $location

    if ($valueIsText) {
        $ws.Cells.Item($row, 5).NumberFormat = "@"
    }
    $ws.Cells.Item($row, 5).Value = $value

    $ws.Cells.Item($row, 6).NumberFormat = "@"
    $ws.Cells.Item($row, 6).Value = $status
}

# --- mmWave(InBed) sheet: append rows 168-172 ---
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")

Add-LogRow $wsInBed 168 "2026-02-01" "21:22:48" "21:00" "Bedroom" "In Bed" $true "Occupied"
Add-LogRow $wsInBed 169 "2026-02-01" "21:22:50" "21:00" "Bedroom" "In Bed" $true "Occupied"
Add-LogRow $wsInBed 170 "2026-02-01" "21:22:51" "21:00" "Bedroom" "In Bed" $true "Occupied"
Add-LogRow $wsInBed 171 "2026-02-01" "21:22:54" "21:00" "Bedroom" "In Bed" $true "Occupied"
Add-LogRow $wsInBed 172 "2026-02-01" "21:23:13" "21:00" "Bedroom" "Out of Bed" $true "Empty"

# --- mmWave(BR) sheet: append rows 163-165 ---
$wsBR = $wb.Worksheets.Item("mmWave(BR)")

Add-LogRow $wsBR 163 "2026-02-01" "21:22:51" "21:00" "Bedroom" 25 $false "Occupied"
Add-LogRow $wsBR 164 "2026-02-01" "21:22:52" "21:00" "Bedroom" 2 $false "Occupied"
Add-LogRow $wsBR 165 "2026-02-01" "21:22:55" "21:00" "Bedroom" 1 $false "Occupied"

# --- mmWave(HR) sheet: append rows 164-166 ---
$wsHR = $wb.Worksheets.Item("mmWave(HR)")

Add-LogRow $wsHR 164 "2026-02-01" "21:22:50" "21:00" "Bedroom" 73 $false "Occupied"
Add-LogRow $wsHR 165 "2026-02-01" "21:22:52" "21:00" "Bedroom" 50 $false "Occupied"
Add-LogRow $wsHR 166 "2026-02-01" "21:22:55" "21:00" "Bedroom" 49 $false "Occupied"
